$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.020335
$ws.Range("H2").Value = 0.061005
$ws.Range("I2").Value = 0.009804808687698561
$ws.Range("J2").Value = 0.009804808687698559
$ws.Range("M2").Value = 0.06743766666666666
$ws.Range("N2").Value = 0.202313
$ws.Range("O2").Value = 0.004349811883262162
$ws.Range("P2").Value = 0.004349811883262163
$ws.Range("Q2").Value = 0.001371344951666667
$ws.Range("R2").Value = 0.012342104565
$ws.Range("S2").Value = 0.00004264907334286329
$ws.Range("T2").Value = 0.00004264907334286329
$ws.Range("G3").Value = 0.020335
$ws.Range("H3").Value = 0.061005
$ws.Range("I3").Value = 0.009804808687698561
$ws.Range("J3").Value = 0.009804808687698559
$ws.Range("O3").Value = 0.7938207485680675
$ws.Range("P3").Value = 0.7938207485680676
$ws.Range("Q3").Value = 0.2502641735533334
$ws.Range("R3").Value = 2.25237756198
$ws.Range("S3").Value = 0.007783260572035563
$ws.Range("T3").Value = 0.007783260572035563
$ws.Range("G4").Value = 0.020335
$ws.Range("H4").Value = 0.061005
$ws.Range("I4").Value = 0.009804808687698561
$ws.Range("J4").Value = 0.009804808687698559
$ws.Range("O4").Value = 0.2018294395486703
$ws.Range("P4").Value = 0.2018294395486704
$ws.Range("Q4").Value = 0.06362982824333333
$ws.Range("R4").Value = 0.57266845419
$ws.Range("S4").Value = 0.001978899042320134
$ws.Range("T4").Value = 0.001978899042320134
$ws.Range("I5").Value = 0.1486140913768632
$ws.Range("J5").Value = 0.1486140913768632
$ws.Range("M5").Value = 0.06743766666666666
$ws.Range("N5").Value = 0.202313
$ws.Range("O5").Value = 0.004349811883262162
$ws.Range("P5").Value = 0.004349811883262163
$ws.Range("Q5").Value = 0.020785839933
$ws.Range("R5").Value = 0.187072559397
$ws.Range("S5").Value = 0.0006464433406912886
$ws.Range("T5").Value = 0.0006464433406912886
$ws.Range("I6").Value = 0.1486140913768632
$ws.Range("J6").Value = 0.1486140913768632
$ws.Range("O6").Value = 0.7938207485680675
$ws.Range("P6").Value = 0.7938207485680676
$ws.Range("S6").Value = 0.1179729492645448
$ws.Range("T6").Value = 0.1179729492645448
$ws.Range("I7").Value = 0.1486140913768632
$ws.Range("J7").Value = 0.1486140913768632
$ws.Range("O7").Value = 0.2018294395486703
$ws.Range("P7").Value = 0.2018294395486704
$ws.Range("S7").Value = 0.02999469877162719
$ws.Range("T7").Value = 0.02999469877162719
$ws.Range("H8").Value = 5.236273000000001
$ws.Range("I8").Value = 0.8415810999354383
$ws.Range("J8").Value = 0.8415810999354382
$ws.Range("M8").Value = 0.06743766666666666
$ws.Range("N8").Value = 0.202313
$ws.Range("O8").Value = 0.004349811883262162
$ws.Range("P8").Value = 0.004349811883262163
$ws.Range("Q8").Value = 0.1177073443832222
$ws.Range("R8").Value = 1.059366099449
$ws.Range("S8").Value = 0.003660719469228011
$ws.Range("T8").Value = 0.003660719469228011
$ws.Range("H9").Value = 5.236273000000001
$ws.Range("I9").Value = 0.8415810999354383
$ws.Range("J9").Value = 0.8415810999354382
$ws.Range("O9").Value = 0.7938207485680675
$ws.Range("P9").Value = 0.7938207485680676
$ws.Range("S9").Value = 0.6680645387314872
$ws.Range("T9").Value = 0.6680645387314873
$ws.Range("H10").Value = 5.236273000000001
$ws.Range("I10").Value = 0.8415810999354383
$ws.Range("J10").Value = 0.8415810999354382
$ws.Range("O10").Value = 0.2018294395486703
$ws.Range("P10").Value = 0.2018294395486704
$ws.Range("Q10").Value = 5.461571209330445
$ws.Range("R10").Value = 49.15414088397401
$ws.Range("S10").Value = 0.169855841734723
$ws.Range("T10").Value = 0.169855841734723
